$d = $word.ActiveDocument

# --- Title paragraph: "The" " " "product" " " "rule"  ->  "The product rule" ---
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.Find.Execute("The product rule", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "The product rule", 2)

# --- Author paragraph: "Tom" " " "Coleman" -> "Tom Coleman" ---
$authorRange = $d.Paragraphs.Item(2).Range
$authorRange.Find.Execute("Tom Coleman", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "Tom Coleman", 2)

# --- Abstract paragraph: many single-word runs -> one run with the full sentence ---
$abstractText = "The product rule is one of the three central techniques of differentiation, " + `
    "allowing you to differentiate a product of two functions. This guide introduces the " + `
    "product rule and explains examples of where it is used."
$abstractRange = $d.Paragraphs.Item(4).Range
$abstractRange.Find.Execute($abstractText, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $abstractText, 2)
